$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 8022
$ws1.Range("F15").Value = 7
$ws1.Range("F20").Value = 1301
$ws1.Range("F24").Value = 1292
$ws1.Range("F35").Value = 149
$ws1.Range("F40").Value = 101

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F32").Value = 1014

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 869
$ws3.Range("F7").Value = 224
$ws3.Range("F10").Value = 2735

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 869
$ws4.Range("F6").Value = 8022
$ws4.Range("F7").Value = 224
$ws4.Range("F12").Value = 2735
$ws4.Range("F20").Value = 1301
$ws4.Range("F24").Value = 1292
$ws4.Range("F36").Value = 149
$ws4.Range("F40").Value = 101
